$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values for A2:A267 (Nivel_Felicidad), reflecting the data-cleaning
# pass applied to the column in the commit.
$values = @(2,2,2,2,2,2,3,2,2,2,3,2,2,2,2,2,2,2,2,2,3,2,2,2,2,2,2,2,2,3,2,3,2,3,3,2,3,2,2,2,2,3,2,2,2,3,2,3,2,2,2,2,2,2,3,3,2,2,2,3,1,2,2,2,2,2,2,2,3,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,3,2,2,2,2,2,2,2,2,2,2,3,2,2,2,2,2,2,2,2,3,2,2,2,2,2,2,2,3,2,2,2,2,2,2,2,2,2,2,2,2,3,2,2,3,2,2,2,2,2,2,2,2,3,3,2,2,2,2,2,3,2,3,2,2,2,2,2,2,2,2,2,2,3,2,2,2,2,2,2,2,2,2,2,2,2,3,2,2,3,2,2,2,3,2,3,2,2,2,3,3,3,2,2,2,2,2,3,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,3,2,2,2,3,2,2,2,3,3,2,2,2,2,2,3,2,3,2,2,2,2,2,2,3,2,2,3,2,2,2,2,2,2,3,2,2,2,2,3,2,2,2,2,2,2,2)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}
